$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 18.68421
$ws.Range("I2").Value = 18.68421
$ws.Range("K2").Value = 18.68421
$ws.Range("M2").Value = 94.31578999999999
$ws.Range("H40").Value = 16668616
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 16668616
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 16668616
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -16668966
$ws.Range("H64").Value = 166673730
$ws.Range("J64").Value = 7899.5
$ws.Range("L64").Value = 7899.5
$ws.Range("N64").Value = -8395.5
$ws.Range("H67").Value = 166673730
$ws.Range("J67").Value = 7899.5
$ws.Range("L67").Value = 7899.5
$ws.Range("N67").Value = -9615.5
$ws.Range("H74").Value = 250005250
$ws.Range("I74").Value = 250005250
$ws.Range("K74").Value = 250005250
$ws.Range("M74").Value = -250004314
$ws.Range("H77").Value = 250005250
$ws.Range("I77").Value = 250005250
$ws.Range("K77").Value = 1250026250
$ws.Range("M77").Value = -1250021570
$ws.Range("H112").Value = 3799.279
$ws.Range("J112").Value = 3799.279
$ws.Range("L112").Value = 11397.837
$ws.Range("N112").Value = -13613.837
$ws.Range("H113").Value = 131592870
$ws.Range("J113").Value = 107161000
$ws.Range("L113").Value = 107161000
$ws.Range("N113").Value = -107167508
$ws.Range("H135").Value = 4546412.5
$ws.Range("I135").Value = 5263793.5
$ws.Range("K135").Value = 47374141.5
$ws.Range("M135").Value = -47371606.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 5533.4287
$ws.Range("I57").Value = 5533.4287
$ws.Range("K57").Value = 5533.4287
$ws.Range("M57").Value = -5049.4287
$ws.Range("H63").Value = 2519
$ws.Range("I63").Value = 2198.6667
$ws.Range("K63").Value = 2198.6667
$ws.Range("M63").Value = -1512.6667
$ws.Range("H66").Value = 2519
$ws.Range("I66").Value = 2198.6667
$ws.Range("K66").Value = 10993.3335
$ws.Range("M66").Value = -7561.333500000001
$ws.Range("H76").Value = 39999.4
$ws.Range("J76").Value = 39999.4
$ws.Range("L76").Value = 39999.4
$ws.Range("N76").Value = -40675.4
$ws.Range("H79").Value = 39999.4
$ws.Range("J79").Value = 39999.4
$ws.Range("L79").Value = 39999.4
$ws.Range("N79").Value = -42339.4
$ws.Range("H122").Value = 5686.8887
$ws.Range("J122").Value = 5014
$ws.Range("L122").Value = 15042
$ws.Range("N122").Value = -19942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7581213
$ws.Range("I20").Value = 11115209
$ws.Range("J20").Value = 8364.286
$ws.Range("K20").Value = 11115209
$ws.Range("L20").Value = 8364.286
$ws.Range("M20").Value = -11114962
$ws.Range("N20").Value = -8858.286
$ws.Range("H82").Value = 16692
$ws.Range("I82").Value = 16692
$ws.Range("K82").Value = 16692
$ws.Range("M82").Value = -16309
$ws.Range("H85").Value = 16692
$ws.Range("I85").Value = 16692
$ws.Range("K85").Value = 16692
$ws.Range("M85").Value = -15366
$ws.Range("H134").Value = 4103265.2
$ws.Range("I134").Value = 5816754
$ws.Range("J134").Value = 9931.444
$ws.Range("K134").Value = 17450262
$ws.Range("L134").Value = 29794.332
$ws.Range("M134").Value = -17447727
$ws.Range("N134").Value = -34864.33199999999
$ws.Range("H140").Value = 34888
$ws.Range("J140").Value = 34888
$ws.Range("L140").Value = 34888
$ws.Range("N140").Value = -45248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 231.93333
$ws.Range("I5").Value = 273.4
$ws.Range("K5").Value = 273.4
$ws.Range("M5").Value = -161.4
$ws.Range("H132").Value = 6790.968
$ws.Range("I132").Value = 4226.25
$ws.Range("K132").Value = 12678.75
$ws.Range("M132").Value = -10148.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1945.9
$ws.Range("J5").Value = 3517.375
$ws.Range("L5").Value = 10552.125
$ws.Range("N5").Value = -10776.125
$ws.Range("H57").Value = 6666
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 6666
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19998
$ws.Range("M57").Value = $null
$ws.Range("N57").Value = -21116
$ws.Range("H132").Value = 8108.341
$ws.Range("I132").Value = 5684
$ws.Range("K132").Value = 51156
$ws.Range("M132").Value = -48626
$ws.Range("H135").Value = 1945.9
$ws.Range("J135").Value = 3517.375
$ws.Range("L135").Value = 31656.375
$ws.Range("N135").Value = -36726.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 42418
$ws.Range("J75").Value = 42418
$ws.Range("L75").Value = 42418
$ws.Range("N75").Value = -44166
$ws.Range("H78").Value = 42418
$ws.Range("J78").Value = 42418
$ws.Range("L78").Value = 127254
$ws.Range("N78").Value = -135990
$ws.Range("H80").Value = 3796
$ws.Range("J80").Value = 3993.3333
$ws.Range("L80").Value = 3993.3333
$ws.Range("N80").Value = -5989.3333
$ws.Range("H83").Value = 3796
$ws.Range("J83").Value = 3993.3333
$ws.Range("L83").Value = 19966.6665
$ws.Range("N83").Value = -29950.6665
$ws.Range("H102").Value = 4804.136
$ws.Range("I102").Value = 4637.35
$ws.Range("K102").Value = 4637.35
$ws.Range("M102").Value = -3015.35
$ws.Range("H113").Value = 6712.2354
$ws.Range("I113").Value = 4345.3335
$ws.Range("J113").Value = 9375
$ws.Range("K113").Value = 4345.3335
$ws.Range("L113").Value = 9375
$ws.Range("M113").Value = -2175.3335
$ws.Range("N113").Value = -13715
$ws.Range("H122").Value = 9055852
$ws.Range("I122").Value = 9055852
$ws.Range("K122").Value = 27167556
$ws.Range("M122").Value = -27165106
$ws.Range("H126").Value = 9317.322
$ws.Range("I126").Value = 4998.75
$ws.Range("K126").Value = 14996.25
$ws.Range("M126").Value = -12526.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3750
$ws.Range("I22").Value = 3750
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3750
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3455
$ws.Range("N22").Value = $null
$ws.Range("H27").Value = 3750
$ws.Range("I27").Value = 3750
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3750
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3643
$ws.Range("N27").Value = $null
$ws.Range("H46").Value = 1337.4117
$ws.Range("I46").Value = 1168.5385
$ws.Range("J46").Value = 1886.25
$ws.Range("K46").Value = 1168.5385
$ws.Range("L46").Value = 1886.25
$ws.Range("M46").Value = -980.5385000000001
$ws.Range("N46").Value = -2262.25
$ws.Range("H61").Value = 4276.6
$ws.Range("I61").Value = 2481
$ws.Range("J61").Value = 8466.333000000001
$ws.Range("K61").Value = 2481
$ws.Range("L61").Value = 8466.333000000001
$ws.Range("M61").Value = -2279
$ws.Range("N61").Value = -8870.333000000001
$ws.Range("H100").Value = 6520.524
$ws.Range("I100").Value = 5436.6665
$ws.Range("K100").Value = 5436.6665
$ws.Range("M100").Value = -4895.6665
$ws.Range("H113").Value = 4276.6
$ws.Range("I113").Value = 2481
$ws.Range("J113").Value = 8466.333000000001
$ws.Range("K113").Value = 2481
$ws.Range("L113").Value = 8466.333000000001
$ws.Range("M113").Value = -311
$ws.Range("N113").Value = -12806.333
$ws.Range("H122").Value = 4521.4185
$ws.Range("I122").Value = 3851.8518
$ws.Range("J122").Value = 5651.3125
$ws.Range("K122").Value = 11555.5554
$ws.Range("L122").Value = 16953.9375
$ws.Range("M122").Value = -9105.555399999999
$ws.Range("N122").Value = -21853.9375
$ws.Range("H136").Value = 11088.429
$ws.Range("I136").Value = 3791.5833
$ws.Range("K136").Value = 11374.7499
$ws.Range("M136").Value = -8824.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 285000
$ws.Range("J62").Value = 70000
$ws.Range("L62").Value = 70000
$ws.Range("N62").Value = -71248
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H65").Value = 285000
$ws.Range("J65").Value = 70000
$ws.Range("L65").Value = 350000
$ws.Range("N65").Value = -356240
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H81").Value = 1559386.6
$ws.Range("I81").Value = 2725671.2
$ws.Range("K81").Value = 5451342.4
$ws.Range("M81").Value = -5450281.4
$ws.Range("H84").Value = 1559386.6
$ws.Range("I84").Value = 2725671.2
$ws.Range("K84").Value = 27256712
$ws.Range("M84").Value = -27251408
$ws.Range("H122").Value = 146688.39
$ws.Range("I122").Value = 237869.17
$ws.Range("J122").Value = 5772.636
$ws.Range("K122").Value = 713607.51
$ws.Range("L122").Value = 17317.908
$ws.Range("M122").Value = -711157.51
$ws.Range("N122").Value = -22217.908
$ws.Range("H126").Value = 1792.2222
$ws.Range("I126").Value = 1376.6666
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4129.9998
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1659.9998
$ws.Range("N126").Value = -10940
